# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with the latest scraped values (GitHub Actions cron update).
#
# Note: several Price values are plain decimals (e.g. "21.84", "1.002")
# that Excel would otherwise auto-convert to numeric cells. They are
# entered with a leading apostrophe so Excel keeps them as literal text,
# matching how the sheet already stores every other Price value (values
# such as "27.780.43" that contain two dots are never auto-numified, so
# no apostrophe is required there).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.780.43"
$ws.Range("E2").Value = "  +5.95%  "
$ws.Range("D3").Value = "1.731.52"
$ws.Range("E3").Value = "  +4.33%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'227.20"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("E6").Value = "  +3.14%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'0.2730"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("D9").Value = "'0.06674"
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("D10").Value = "'21.84"
$ws.Range("E10").Value = "  +5.70%  "
$ws.Range("D11").Value = "'0.07784"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "'4.685"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").Value = "1.725.77"
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").Value = "1.969.48"
$ws.Range("E14").Value = "  +4.34%  "
$ws.Range("D15").Value = "'0.5954"
$ws.Range("E15").Value = "  +5.42%  "
$ws.Range("D16").Value = "0.0₅8390"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "'68.96"
$ws.Range("E17").Value = "  +4.80%  "
$ws.Range("D18").Value = "27.772.99"
$ws.Range("E18").Value = "  +5.98%  "
$ws.Range("D19").Value = "'226.21"
$ws.Range("E19").Value = "  +17.91%  "
$ws.Range("D20").Value = "'4.803"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +4.10%  "
$ws.Range("D23").Value = "'6.202"
$ws.Range("E23").Value = "  +3.31%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "'147.12"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "'1.720"
$ws.Range("E26").Value = "  +12.58%  "
$ws.Range("D27").Value = "'0.1249"
$ws.Range("E27").Value = "  +3.69%  "
$ws.Range("D28").Value = "'7.460"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").Value = "'17.04"
$ws.Range("E29").Value = "  +6.00%  "
$ws.Range("D30").Value = "'0.05663"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D32").Value = "'3.655"
$ws.Range("E32").Value = "  +4.44%  "
$ws.Range("D33").Value = "'3.502"
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("E34").Value = "  +5.62%  "
$ws.Range("D35").Value = "'0.9731"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "'0.5975"
$ws.Range("E38").Value = "  +3.52%  "
$ws.Range("E39").Value = "  +4.09%  "
$ws.Range("D40").Value = "'5.908"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "'0.8604"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "1.048.44"
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").Value = "'101.56"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Value = "1.874.08"
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("E46").Value = "  +8.96%  "
$ws.Range("D47").Value = "'59.58"
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").Value = "'8.262"
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").Value = "'0.4429"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").Value = "'0.05332"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "'0.9995"
$ws.Range("E51").Value = "  -0.56%  "
